# Apply updated "想去人数" (F) and "最低票价" (G) values to the
# "展览" and "全部类型" worksheets, matching the upstream data refresh.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1581
$ws1.Range("G2").Value = 65
$ws1.Range("F3").Value = 244
$ws1.Range("G3").Value = 75
$ws1.Range("F4").Value = 108
$ws1.Range("F7").Value = 6446
$ws1.Range("F9").Value = 418
$ws1.Range("F10").Value = 128
$ws1.Range("F11").Value = 5739
$ws1.Range("F17").Value = 76
$ws1.Range("F18").Value = 376
$ws1.Range("F21").Value = 325
$ws1.Range("F24").Value = 4091
$ws1.Range("F25").Value = 24

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1581
$ws4.Range("G2").Value = 65
$ws4.Range("F3").Value = 244
$ws4.Range("G3").Value = 75
$ws4.Range("F4").Value = 108
$ws4.Range("F7").Value = 6446
$ws4.Range("F9").Value = 418
$ws4.Range("F10").Value = 128
$ws4.Range("F11").Value = 5739
$ws4.Range("F17").Value = 76
$ws4.Range("F18").Value = 376
$ws4.Range("F21").Value = 325
$ws4.Range("F24").Value = 4091
$ws4.Range("F26").Value = 24
